$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the four "FLD_Transmittals_New_*" tests as Y for Sanity Runmode (D5:D8)
$ws.Range("D5:D8").Value = "Y"

# Add new test row for FLD_Transmittals_LeftNavigationBar
$ws.Cells.Item(15, 1).Value = "FLD_Transmittals_LeftNavigationBar"
$ws.Cells.Item(15, 2).Value = "Verifies the Left Navigation menu items"
$ws.Cells.Item(15, 3).Value = "N"
$ws.Cells.Item(15, 4).Value = "Y"
$ws.Cells.Item(15, 6).Value = "Sprint2"

# Copy formatting of the last existing row onto the new row
$ws.Range("A14:F14").Copy()
$ws.Range("A15:F15").PasteSpecial(-4122)

# Extend the data validation ranges to include the new row
$ws.Range("C2:D15").Validation.Delete()
$ws.Range("C2:D15").Validation.Add(3, 1, 1, '"Y,N"')
$ws.Range("F2:F15").Validation.Delete()
$ws.Range("F2:F15").Validation.Add(3, 1, 1, '"Sprint1,Sprint2,Sprint3,Sprint4,Sprint5,Sprint6,Sprint7,Sprint8,Sprint9,Sprint10"')

# Restore last active selection
$ws.Range("D5").Select()
